$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5
$ws.Range("C4").Value = 1.4

# Update column widths
$ws.Columns.Item(1).ColumnWidth = 27
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(3).ColumnWidth = 27.25

# Update selection
$ws.Range("F4").Select()
